# Apply updates to RPAR_holdings workbook:
#  1) Update the "as of" date in the confidential disclosure note (A18)
#     from 2021-07-08 to 2021-07-09.
#  2) Update the Weight (D) and Percent Change (E) values for rows 2-15
#     on Sheet1 to the refreshed figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Sheet is protected; unprotect to allow edits, then restore protection after.
$ws.Unprotect()

# --- 1) Update disclosure text in A18 ---
$oldText = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution." + [char]10 + "Model holdings provided as of 2021-07-08 for illustrative purposes only and are subject to change."
$newText = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution." + [char]10 + "Model holdings provided as of 2021-07-09 for illustrative purposes only and are subject to change."

$cellA18 = $ws.Range("A18")
if ($cellA18.Value -eq $oldText) {
    $cellA18.Value = $newText
} else {
    $cellA18.Value = $newText
}

# --- 2) Update Weight (D) / Percent Change (E) values for rows 2-14 ---
$updates = @(
    @{ Row = 2;  D = 0.05821814849958101;  E = 0.01209785823102427 },
    @{ Row = 3;  D = 0.01994776331106263;  E = 0.01726167124362488 },
    @{ Row = 4;  D = 0.0271448182313218;   E = 0.01825869690563131 },
    @{ Row = 5;  D = 0.02851431943941751;  E = 0.01905132192846026 },
    @{ Row = 6;  D = 0.02892256122811651;  E = 0.03720405862457721 },
    @{ Row = 7;  D = 0.01763917556027273;  E = 0.02148033126293969 },
    @{ Row = 8;  D = 0.01057254932649887;  E = 0.002590673575129543 },
    @{ Row = 9;  D = 0.0105412464417424;   E = 0.01373422420193005 },
    @{ Row = 10; D = 0.06653819441047865;  E = 0.003351955307262733 },
    @{ Row = 11; D = 0.06664971093742358;  E = 0.003346346904628872 },
    @{ Row = 12; D = 0.1550744910835689;   E = -0.01406271026779704 },
    @{ Row = 13; D = 0.3967640642882996;   E = -0.005042735042735114 },
    @{ Row = 14; D = 0.1134729572422159;   E = -0.00295517241379295 }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 4).Value = $u.D
    $ws.Cells.Item($u.Row, 5).Value = $u.E
}

# Row 15 only has Percent Change (E) updated; D15 (Total weight = 1) stays the same.
$ws.Cells.Item(15, 5).Value = -0.0003562007427913016

# Restore sheet protection to match original protected state.
$ws.Protect()
